# Add data for 2022-11-24
# Updates the 2022 year-to-date column (and a couple of adjacent 2021 cells
# that shifted week) across the citywide summary, neighborhood roll-up, and
# each per-neighborhood sheet.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("I2").Value = 6586
$ws.Range("I3").Value = 6877
$ws.Range("H4").Value = 1676
$ws.Range("I4").Value = 1577
$ws.Range("I5").Value = 639
$ws.Range("I6").Value = 7895
$ws.Range("H7").Value = 25988
$ws.Range("I7").Value = 23574

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("I3").Value = 63
$ws.Range("I7").Value = 265

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("I3").Value = 35
$ws.Range("I7").Value = 131

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("I3").Value = 28
$ws.Range("I7").Value = 79

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("I2").Value = 116
$ws.Range("I7").Value = 417

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("I2").Value = 215
$ws.Range("I3").Value = 335
$ws.Range("I7").Value = 906

$ws = $wb.Worksheets.Item("New City")
$ws.Range("I2").Value = 181
$ws.Range("I6").Value = 160
$ws.Range("I7").Value = 541

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("I2").Value = 184
$ws.Range("I5").Value = 71
$ws.Range("I7").Value = 746
$ws.Range("I8").Value = 1409
$ws.Range("I11").Value = 358
$ws.Range("I14").Value = 131
$ws.Range("I15").Value = 275
$ws.Range("I19").Value = 667
$ws.Range("I20").Value = 586
$ws.Range("I21").Value = 103
$ws.Range("I22").Value = 67
$ws.Range("I29").Value = 1427
$ws.Range("I30").Value = 79
$ws.Range("I33").Value = 1058
$ws.Range("I36").Value = 324
$ws.Range("I44").Value = 175
$ws.Range("I47").Value = 170
$ws.Range("I51").Value = 281
$ws.Range("I53").Value = 258
$ws.Range("I54").Value = 476
$ws.Range("I59").Value = 40
$ws.Range("H63").Value = 227
$ws.Range("I63").Value = 74
$ws.Range("I65").Value = 541
$ws.Range("I66").Value = 65
$ws.Range("I67").Value = 906
$ws.Range("I73").Value = 217
$ws.Range("I75").Value = 74
$ws.Range("I76").Value = 341
$ws.Range("I79").Value = 669
$ws.Range("I83").Value = 511
$ws.Range("I86").Value = 149
$ws.Range("I90").Value = 305
$ws.Range("I94").Value = 239
$ws.Range("I95").Value = 356
$ws.Range("I96").Value = 265
$ws.Range("I99").Value = 417
$ws.Range("I100").Value = 40
$ws.Range("H101").Value = 25988
$ws.Range("I101").Value = 23574

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("I2").Value = 171
$ws.Range("I7").Value = 511

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("I2").Value = 124
$ws.Range("I7").Value = 356

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("I3").Value = 389
$ws.Range("I6").Value = 340
$ws.Range("I7").Value = 1058

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("I3").Value = 102
$ws.Range("I7").Value = 476

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("I2").Value = 419
$ws.Range("I3").Value = 494
$ws.Range("I4").Value = 73
$ws.Range("I7").Value = 1427

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("I3").Value = 194
$ws.Range("I5").Value = 18
$ws.Range("I7").Value = 667

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("I2").Value = 55
$ws.Range("I7").Value = 175

$ws = $wb.Worksheets.Item("River North")
$ws.Range("I2").Value = 69
$ws.Range("I7").Value = 341

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("I3").Value = 259
$ws.Range("I4").Value = 56

$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("I6").Value = 78
$ws.Range("I7").Value = 103

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("I5").Value = 26
$ws.Range("I7").Value = 669

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("I2").Value = 165
$ws.Range("I7").Value = 586

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("I3").Value = 108
$ws.Range("I6").Value = 102
$ws.Range("I7").Value = 324

$ws = $wb.Worksheets.Item("Wrigleyville")
$ws.Range("I2").Value = 7
$ws.Range("I6").Value = 40

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("I3").Value = 40
$ws.Range("I7").Value = 239

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("I6").Value = 56
$ws.Range("I7").Value = 170

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("I2").Value = 80
$ws.Range("I7").Value = 275

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("I3").Value = 12
$ws.Range("I7").Value = 65

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("I3").Value = 77
$ws.Range("I7").Value = 358

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("I4").Value = 22
$ws.Range("I6").Value = 57
$ws.Range("I7").Value = 217

$ws = $wb.Worksheets.Item("Montclare")
$ws.Range("I4").Value = 4
$ws.Range("I7").Value = 40

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("I6").Value = 38
$ws.Range("I7").Value = 184

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("I2").Value = 418
$ws.Range("I3").Value = 407
$ws.Range("I7").Value = 1409

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("I2").Value = 20
$ws.Range("I7").Value = 71

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("I2").Value = 28
$ws.Range("I7").Value = 149

$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("I2").Value = 23
$ws.Range("I7").Value = 74

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("I3").Value = 76
$ws.Range("I7").Value = 305

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("I6").Value = 111
$ws.Range("I7").Value = 281

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("I2").Value = 57
$ws.Range("I7").Value = 258

$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("I6").Value = 18
$ws.Range("I7").Value = 67

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("I3").Value = 230
$ws.Range("I7").Value = 746

Write-Output "Updated 131 cells across 39 sheets"
